$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.046.65"
$ws.Range("E2").Value = '  +2.46%  '
$ws.Range("D3").Value = "'2.359.47"
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = "'311.31"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = "'107.57"
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = '  -2.38%  '
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("D11").Value = "'0.0914"
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").Value = "'0.971"
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").Value = "'2.718.92"
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = "'15.15"
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = "'2.364.88"
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = "'45.012.40"
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = "'14.19"
$ws.Range("E19").Value = '  +9.35%  '
$ws.Range("E20").Value = '  -4.79%  '
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = "'72.77"
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("D23").Value = "'3.49"
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("D24").Value = "'257.82"
$ws.Range("E24").Value = '  -4.02%  '
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = "'11.04"
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("D28").Value = "'7.16"
$ws.Range("E28").Value = '  -5.95%  '
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("E30").Value = '  +8.70%  '
$ws.Range("D31").Value = "'22.24"
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").Value = "'37.09"
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("D33").Value = "'167.76"
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("E34").Value = '  +4.52%  '
$ws.Range("E36").Value = '  +1.14%  '
$ws.Range("D37").Value = "'4.67"
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = '  +3.64%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = "'2.90"
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.0351"
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("D42").Value = "'99.65"
$ws.Range("E42").Value = '  -4.51%  '
$ws.Range("D43").Value = "'1.879.83"
$ws.Range("E43").Value = '  +12.21%  '
$ws.Range("D44").Value = "'69.16"
$ws.Range("E44").Value = '  -3.29%  '
$ws.Range("E45").Value = '  -4.06%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").Value = "'12.77"
$ws.Range("E47").Value = '  -4.73%  '
$ws.Range("D48").Value = "'81.10"
$ws.Range("E48").Value = '  +5.38%  '
$ws.Range("E49").Value = '  +8.16%  '
$ws.Range("D50").Value = "'110.07"
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("E51").Value = '  +2.28%  '
